$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing CP-like values in column E for the first two data rows
$ws.Range("E2").Value = 3000
$ws.Range("E3").Value = 1000

# Add two new rows (4 and 5) that replicate rows 2 and 3 respectively
# (same idCliente/Operativa/DOMICILIO/LOCALIDAD/pesoenvio), but with a
# different quantity so the sheet now covers any branch/sucursal.
$ws.Range("A2:F2").Copy($ws.Range("A4:F4"))
$ws.Range("A3:F3").Copy($ws.Range("A5:F5"))

$ws.Range("E4").Value = 100000
$ws.Range("E5").Value = 5000

# Move the active selection as recorded in the saved workbook
$ws.Range("E6").Select() | Out-Null

$wb.Save()
